$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7: Offices/Manufacturing weights dropped to 0, weight count 13 -> 12 ---
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0
$ws.Range("AC7").Value = 12

# --- LockdownEffectiveness (column AC) recalculated for rows 24-221 ---
$ws.Range("AC24").Value = 0.04166666666666666
$ws.Range("AC25:AC26").Value = 0.2499999999916667
$ws.Range("AC27").Value = 0.4791666666499999
$ws.Range("AC28:AC30").Value = 0.5624999999833332
$ws.Range("AC31:AC33").Value = 0.6458333333166666
$ws.Range("AC34:AC65").Value = 0.8333333333166667
$ws.Range("AC66:AC70").Value = 0.7499999999833333
$ws.Range("AC71:AC81").Value = 0.5833333333166665
$ws.Range("AC82:AC96").Value = 0.5624999999833332
$ws.Range("AC97:AC104").Value = 0.4791666666583334
$ws.Range("AC105:AC118").Value = 0.3541666666583334
$ws.Range("AC119:AC129").Value = 0.2083333333333333
$ws.Range("AC130:AC146").Value = 0.2361111111083334
$ws.Range("AC147:AC221").Value = 0.3194444444416667

# --- Append 12 new daily rows (9/30/2020 .. 10/11/2020), rows 222-233 ---
$newDates = @("9/30/2020","10/1/2020","10/2/2020","10/3/2020","10/4/2020","10/5/2020","10/6/2020","10/7/2020","10/8/2020","10/9/2020","10/10/2020","10/11/2020")
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")
$rowVals = @(0,0,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,0)

# Copy the formatting (bold/border/alignment style) of the last date cell down
$ws.Range("A221").Copy() | Out-Null
$ws.Range("A222:A233").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = 222 + $i
    $addr = "A" + $r

    # Write the date text as a shared-string (not an auto-converted date
    # serial) by routing it through a text formula, then freezing the
    # result back down to a plain value via copy / paste-values.
    $ws.Range($addr).Formula = "=""" + $newDates[$i] + """"
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false

    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $r).Value = $rowVals[$j]
    }

    $ws.Range("AC" + $r).Value = 0.3194444444416667
}
